# "Generate Report for Handback" — refresh the handoff/handback timestamps
# for the 9b702496-0af6-4b72-b27c-ad52d2bb0462.md file across the Overview,
# zh-cn, and de-de sheets, as if a new handback report had just been
# generated for that file.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the .md file (row 4)
$wsOverview.Range("G4").Value = "2016-11-29 04:26:43"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the same .md file (row 4)
$wsZhCn.Range("H4").Value = "2016-11-29 04:26:28"
$wsZhCn.Range("K4").Value = "2016-11-29 04:27:20"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the same .md file (row 4)
$wsDeDe.Range("H4").Value = "2016-11-29 04:26:43"
$wsDeDe.Range("K4").Value = "2016-11-29 04:27:38"
